$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 14:41"

# --- Countries reordered: swap labels so the data stays attached to the
#     correct row index mapping used by the sharedStrings table. ---
# Marruecos / Paises Bajos swap (rows 33-34)
$ws.Range("A33").Value = "Paises Bajos"
$ws.Range("A34").Value = "Marruecos"

# Bonaire, San Eustaquio y Saba / Liechtenstein swap (rows 195-196)
$ws.Range("A195").Value = "Liechtenstein"
$ws.Range("A196").Value = "Bonaire, San Eustaquio y Saba"

# --- Updated numeric stats per country row ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7553072
$ws.Range("C4").Value = 3749
$ws.Range("D4").Value = 4777579
$ws.Range("E4").Value = 2561950
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 213543

# Row 18: Irak
$ws.Range("B18").Value = 375931
$ws.Range("C18").Value = 3672
$ws.Range("D18").Value = 303665
$ws.Range("E18").Value = 62919
$ws.Range("G18").Value = 49
$ws.Range("H18").Value = 9347

# Row 33: now Paises Bajos
$ws.Range("B33").Value = 131889
$ws.Range("C33").Value = 3967
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("G33").Value = 21
$ws.Range("H33").Value = 6449

# Row 34: now Marruecos
$ws.Range("B34").Value = 128565
$ws.Range("D34").Value = 106044
$ws.Range("E34").Value = 20258
$ws.Range("H34").Value = 2263

# Row 35: Catar
$ws.Range("B35").Value = 126339
$ws.Range("C35").Value = 175
$ws.Range("D35").Value = 123302
$ws.Range("E35").Value = 2821
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 216

# Row 40: Kuwait
$ws.Range("B40").Value = 106458
$ws.Range("C40").Value = 371
$ws.Range("D40").Value = 98435
$ws.Range("E40").Value = 7403
$ws.Range("G40").Value = 5
$ws.Range("H40").Value = 620

# Row 57: Barein
$ws.Range("E57").Value = 5362
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 257

# Row 66: Kirguistan
$ws.Range("B66").Value = 47184
$ws.Range("C66").Value = 128
$ws.Range("D66").Value = 43278
$ws.Range("E66").Value = 2840
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 1066

# Row 79: Dinamarca
$ws.Range("B79").Value = 29302
$ws.Range("C79").Value = 370
$ws.Range("D79").Value = 22297
$ws.Range("E79").Value = 6351
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 654

# Row 90: Madagascar
$ws.Range("B90").Value = 16529
$ws.Range("C90").Value = 36
$ws.Range("D90").Value = 15480
$ws.Range("E90").Value = 817

# Row 143: Sri Lanka
$ws.Range("D143").Value = 3254
$ws.Range("E143").Value = 121

# Row 195: now Liechtenstein
$ws.Range("B195").Value = 123
$ws.Range("C195").Value = 3
$ws.Range("D195").Value = 116
$ws.Range("E195").Value = 6

# Row 196: now Bonaire, San Eustaquio y Saba
$ws.Range("B196").Value = 121
$ws.Range("D196").Value = 32
$ws.Range("E196").Value = 88
